$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample Section")

# 1. Remove the existing cell comments first. Deleting a column does not
#    re-anchor legacy cell comments to their new column, so the old A1:V1
#    comments must be removed and rebuilt after the shift below.
foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")) {
    $cmt = $ws.Range($col + "1").Comment
    if ($cmt -ne $null) { $cmt.Delete() }
}

# 2. Delete column A ("header_info") - shifts B:V left to A:U
$ws.Range("A1").EntireColumn.Delete()

# 3. Re-create the header comments shifted one column to the left (the old
#    A1 comment was empty and is simply dropped with its column).
$c_A = @'
(Required) Unique identifier for the source (parent) from which the sample was
taken. An example value would be HBM122.EFGH.789.
'@
$ws.Range("A1").AddComment($c_A)

$c_B = @'
(Required) The HuBMAP ID for the sample assigned by the ingest portal. Example:
HBM743.CKJW.876
'@
$ws.Range("B1").AddComment($c_B)

$c_C = @'
An internal field labs can use it to add whatever ID(s) they want or need for
dataset validation and tracking. This could be a single ID (e.g.,
"Visium_9OLC_A4_S1") or a delimited list of IDs (e.g., “9OL; 9OLC.A2;
Visium_9OLC_A4_S1”). This field will not be accessible to anyone outside of the
consortium and no effort will be made to check if IDs provided by one data
provider are also used by another.
'@
$ws.Range("C1").AddComment($c_C)

$c_D = @'
(Required) How long was the source material (parent) stored, prior to this
sample being processed.
'@
$ws.Range("D1").AddComment($c_D)

$c_E = @'
(Required) The time duration unit of measurement
'@
$ws.Range("E1").AddComment($c_E)

$c_F = @'
(Required) DOI for the protocols.io page that describes the assay or sample
procurement and preparation. For example for an imaging assay, the protocol
might include staining of a section through the creation of an OME-TIFF file. In
this case the protocol would include any image processing steps required to
create the OME-TIFF file.
'@
$ws.Range("F1").AddComment($c_F)

$c_G = @'
(Required) The medium used during the sample preparation
'@
$ws.Range("G1").AddComment($c_G)

$c_H = @'
(Required) The condition by which the preparation occurred, such as was the
sample placed in dry ice during the preparation.
'@
$ws.Range("H1").AddComment($c_H)

$c_I = @'
How long the tissue was being handled before the initial preservation.
'@
$ws.Range("I1").AddComment($c_I)

$c_J = @'
The time unit of measurement
'@
$ws.Range("J1").AddComment($c_J)

$c_K = @'
(Required) What was the sample preserved in.
'@
$ws.Range("K1").AddComment($c_K)

$c_L = @'
(Required) The method by which the sample was stored, after preparation and
before the assay was performed.
'@
$ws.Range("L1").AddComment($c_L)

$c_M = @'
For example, RIN: 8.7. For suspensions, measured by visual inspection prior to
cell lysis or defined by known parameters such as wells with several cells or no
cells. This can be captured at a high level. "OK" or "not OK", or with more
specificity such as "debris", "clump", "low clump".
'@
$ws.Range("M1").AddComment($c_M)

$c_N = @'
histopathological reporting of key variables that are important for the tissue
(absence of necrosis, comment on composition, significant pathology description,
high level inflammation/fibrosis assessment, etc.)
'@
$ws.Range("N1").AddComment($c_N)

$c_O = @'
(Required) Thickness of the sample section.
'@
$ws.Range("O1").AddComment($c_O)

$c_P = @'
(Required) The thickness unit of measuement
'@
$ws.Range("P1").AddComment($c_P)

$c_Q = @'
(Required) The index number for the section if the sample is a single section.
'@
$ws.Range("Q1").AddComment($c_Q)

$c_R = @'
The area of the sample section.
'@
$ws.Range("R1").AddComment($c_R)

$c_S = @'
The area unit of measurement
'@
$ws.Range("S1").AddComment($c_S)

$c_T = @'
Miscellaneous details about the sample, not captured in the existing metadata
fields.
'@
$ws.Range("T1").AddComment($c_T)

$c_U = @'
(Required) The string that serves as the definitive identifier for the metadata
schema version and is readily interpretable by computers for data validation and
processing.
'@
$ws.Range("U1").AddComment($c_U)

# 4. Bump the .metadata sheet pav:createdOn timestamp
$meta = $wb.Worksheets.Item(".metadata")
$meta.Range("D2").Value = "2023-10-03T09:51:12-07:00"

Write-Host "done"
